# 2ИСИП-722_ДисМат_.xlsx — grade entry updates (row 10, 19, 20, 26, 27)
# plus the view/selection state that was active when the author last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Кердан Константин) ---
$ws.Range("S10").Value = 4

# --- Row 19 (Петров Иван) ---
$ws.Range("J19").Value = 5
$ws.Range("K19").Value = 10
$ws.Range("S19").Value = 5

# --- Row 20 (Рейзвих Герман) ---
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 10
$ws.Range("S20").Value = 5

# --- Row 26 (Хамуков Тимур) ---
$ws.Range("D26:J26").Value = 5
$ws.Range("K26").Value = 10
$ws.Range("S26").Value = 4

# --- Row 27 (Хасбулатов Магомед) ---
$ws.Range("D27:K27").Value = 5
$ws.Range("S27").Value = 3

# --- View state: scroll the frozen bottom-right pane and land the
# selection on S19, matching the author's last saved cursor position. ---
$window = $excel.ActiveWindow
try {
    $window.ScrollRow = 6
    $window.ScrollColumn = 3
} catch {
    # Older/limited hosts may not expose independent pane scrolling;
    # the explicit cell selection below still records the active cell.
}

$ws.Range("S19").Select()
